$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The import needs a "Route" column in addition to "From"/"To"/"Notes".
# Insert a new column before the existing "Notes" column (C) so "Notes"
# shifts right to D, then give the new C1 cell the "Route" header.
$ws.Columns("C:C").Insert()

# New "Route" header keeps the same look the other headers already use
# (bold text on a yellow fill).
$ws.Range("C1").Value = "Route"
$ws.Range("C1").Font.Bold = $true
$ws.Range("C1").Interior.Color = 65535

# The shifted "Notes" header keeps its yellow fill but is no longer bold.
$ws.Range("D1").Font.Bold = $false
$ws.Range("D1").Interior.Color = 65535
